$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "BBDD - Error Actual" (sheet1): add 5 new columns FU:FY
# Headers: "Usado en E1".."Usado en E5", data = 0/1 flags per row (2-8)
# -------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BBDD - Error Actual")

# Copy header style (bold/border/center) from FT1 into the new header cells
$ws1.Range("FT1").Copy($ws1.Range("FU1:FY1"))

$ws1.Range("FU1").Value = "Usado en E1"
$ws1.Range("FV1").Value = "Usado en E2"
$ws1.Range("FW1").Value = "Usado en E3"
$ws1.Range("FX1").Value = "Usado en E4"
$ws1.Range("FY1").Value = "Usado en E5"

$ws1.Range("FU2").Value = 0
$ws1.Range("FV2").Value = 0
$ws1.Range("FW2").Value = 0
$ws1.Range("FX2").Value = 0
$ws1.Range("FY2").Value = 1

$ws1.Range("FU3").Value = 0
$ws1.Range("FV3").Value = 0
$ws1.Range("FW3").Value = 0
$ws1.Range("FX3").Value = 0
$ws1.Range("FY3").Value = 1

$ws1.Range("FU4").Value = 0
$ws1.Range("FV4").Value = 0
$ws1.Range("FW4").Value = 0
$ws1.Range("FX4").Value = 0
$ws1.Range("FY4").Value = 1

$ws1.Range("FU5").Value = 1
$ws1.Range("FV5").Value = 0
$ws1.Range("FW5").Value = 0
$ws1.Range("FX5").Value = 0
$ws1.Range("FY5").Value = 1

$ws1.Range("FU6").Value = 1
$ws1.Range("FV6").Value = 1
$ws1.Range("FW6").Value = 1
$ws1.Range("FX6").Value = 1
$ws1.Range("FY6").Value = 1

$ws1.Range("FU7").Value = 1
$ws1.Range("FV7").Value = 1
$ws1.Range("FW7").Value = 0
$ws1.Range("FX7").Value = 1
$ws1.Range("FY7").Value = 1

$ws1.Range("FU8").Value = 0
$ws1.Range("FV8").Value = 0
$ws1.Range("FW8").Value = 0
$ws1.Range("FX8").Value = 0
$ws1.Range("FY8").Value = 1

# -------------------------------------------------------------------------
# Sheet "Inicial y final" (sheet3): add 3 new columns P:R
# Headers: "V_sub actual", "Error Actual", "Decaímiento anual"
# -------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Inicial y final")

$ws3.Range("O1").Copy($ws3.Range("P1:R1"))
$ws3.Range("P1").Value = "V_sub actual"
$ws3.Range("Q1").Value = "Error Actual"
$ws3.Range("R1").Value = "Decaímiento anual"

$ws3.Range("P2").Value = 0
$ws3.Range("Q2").Value = 0

$ws3.Range("P3").Value = 0
$ws3.Range("Q3").Value = 0

$ws3.Range("P4").Value = 7.419324184763042
$ws3.Range("Q4").Value = 0.0785784504265742
$ws3.Range("R4").Value = -2.560753770297293

$ws3.Range("P5").Value = 7.419324184763042
$ws3.Range("Q5").Value = 0.1338761216219036
$ws3.Range("R5").Value = -2.558698845263451

# -------------------------------------------------------------------------
# Sheets "RESUMEN E1".."RESUMEN E5" (sheet4..sheet8): add rows 11-13
# Row 11: A = " " (blank separator)
# Row 12: A = "CAPEX", B = <capex value>
# Row 13: A = "Número de medidores", B = <count value>
# -------------------------------------------------------------------------
$resumenSheets = @("RESUMEN E1", "RESUMEN E2", "RESUMEN E3", "RESUMEN E4", "RESUMEN E5")
$capexValues = @(1166756.659902354, 345587.1734157469, 172793.5867078734, 345587.1734157469, 1519338.377610228)
$countValues = @(3, 2, 1, 2, 7)

for ($i = 0; $i -lt $resumenSheets.Length; $i++) {
    $ws = $wb.Worksheets.Item($resumenSheets[$i])

    $ws.Range("A11").Value = " "

    $ws.Range("A12").Value = "CAPEX"
    $ws.Range("B12").Value = $capexValues[$i]

    $ws.Range("A13").Value = "Número de medidores"
    $ws.Range("B13").Value = $countValues[$i]
}
